$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old row 9 data (A9:D9) - values removed while keeping styles.
$ws.Range("A9:D9").ClearContents()

# Update the selected cell to match the new active selection.
$ws.Range("F16").Select()
